$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3259
$wsExhibit.Range("F3").Value = 7
$wsExhibit.Range("F5").Value = 1178
$wsExhibit.Range("F6").Value = 309

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3259
$wsAll.Range("F3").Value = 7
$wsAll.Range("F5").Value = 1178
$wsAll.Range("F7").Value = 309
